$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.252.05"
$ws.Range("E2").Value = "  +0.33%  "
$ws.Range("D3").Value = "1.863.58"
$ws.Range("E3").Value = "  +0.06%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "237.05"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.45%  "
$ws.Range("E6").Value = "  +0.07%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4683"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.50%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2869"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.39%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06552"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.18%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "22.11"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +9.81%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07907"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.13%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "98.02"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.98%  "
$ws.Range("D13").Value = "1.869.94"
$ws.Range("E13").Value = "  +0.82%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.198"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.60%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6866"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.62%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "278.22"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.93%  "
$ws.Range("D17").Value = "30.259.93"
$ws.Range("E17").Value = "  +0.26%  "
$ws.Range("E18").Value = "  +7.86%  "
$ws.Range("E19").Value = "  +0.07%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000007366"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.75%  "
$ws.Range("D21").Value = "2.116.38"
$ws.Range("E21").Value = "  +0.77%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.347"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.88%  "
$ws.Range("E23").Value = "  +0.10%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.204"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.06%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "168.32"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.87%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.264"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.53%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "19.08"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.06%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.963"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +3.06%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.387"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +3.32%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.09862"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.74%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.387"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.27%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.484"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.16%  "
$ws.Range("E33").Value = "  -0.87%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.04739"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.85%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.139"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +4.06%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7049"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.66%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.709"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.04%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01881"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.56%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.626"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +4.35%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.295"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.32%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "75.69"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +4.37%  "
$ws.Range("E42").Value = "  +2.00%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.8524"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.12%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.4185"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.69%  "
$ws.Range("E45").Value = "  +0.05%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "103.58"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.29%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "7.225"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.87%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "954.10"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.42%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.286"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.27%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "34.29"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.42%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05647"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.03%  "
